$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I7").Value = 5
$ws.Range("L7").Value = 5
$ws.Range("X7").Value = 5
$ws.Range("AP7").Value = 98
$ws.Range("AQ7").Value = 249
$ws.Range("AS7").Value = 74
$ws.Range("AT7").Value = 245
$ws.Range("AU7").Value = 0
$ws.Range("AV7").Value = 62
$ws.Range("AW7").Value = 248
$ws.Range("AY7").Value = 52.5
$ws.Range("AZ7").Value = 227
$ws.Range("BB7").Value = 71
$ws.Range("BC7").Value = 230
$ws.Range("BE7").Value = 84
$ws.Range("BF7").Value = 228
$ws.Range("BG7").Value = 0.2199918496694308
$ws.Range("BH7").Value = 0.4840675790259474
$ws.Range("BI7").Value = 0.7584608065647381
$ws.Range("BJ7").Value = 0.1595579165732227
$ws.Range("BK7").Value = 0.3620032855624651
$ws.Range("BL7").Value = 0.8168302750401147
$ws.Range("BN7").Value = 0.1297171271880406
$ws.Range("BO7").Value = 0.3715187419335101
$ws.Range("BQ7").Value = 0.1393776221114097
$ws.Range("BR7").Value = 0.3857636096137051
$ws.Range("BT7").Value = 0.056342058562631
$ws.Range("BU7").Value = 0.5379227053140097
$ws.Range("M8").Value = 10
$ws.Range("R8").Value = 5
$ws.Range("T8").Value = 1
$ws.Range("U8").Value = 6
$ws.Range("Y8").Value = 10
$ws.Range("AC8").Value = 2
$ws.Range("AK8").Value = 9
$ws.Range("AN8").Value = 9
$ws.Range("AO8").Value = 41
$ws.Range("AP8").Value = 137
$ws.Range("AQ8").Value = 216
$ws.Range("AR8").Value = 23
$ws.Range("AS8").Value = 103
$ws.Range("AT8").Value = 219
$ws.Range("AU8").Value = 16
$ws.Range("AV8").Value = 89
$ws.Range("AX8").Value = 16
$ws.Range("AY8").Value = 41
$ws.Range("AZ8").Value = 167
$ws.Range("BA8").Value = 46
$ws.Range("BB8").Value = 75
$ws.Range("BC8").Value = 171
$ws.Range("BD8").Value = 52
$ws.Range("BE8").Value = 86
$ws.Range("BF8").Value = 166
$ws.Range("BG8").Value = 0.3674393175967152
$ws.Range("BH8").Value = 0.5081554580896692
$ws.Range("BI8").Value = 0.6702758497203591
$ws.Range("BJ8").Value = 0.266748786860769
$ws.Range("BK8").Value = 0.4241891200148915
$ws.Range("BL8").Value = 0.836241690408357
$ws.Range("BM8").Value = 0.008546429221175591
$ws.Range("BN8").Value = 0.2361582005224207
$ws.Range("BO8").Value = 0.3353329333524587
$ws.Range("BP8").Value = 0.02056887902668939
$ws.Range("BQ8").Value = 0.1614934391458424
$ws.Range("BR8").Value = 0.3300403890540208
$ws.Range("BS8").Value = 0.001008064516129032
$ws.Range("BT8").Value = 0.08737958077165126
$ws.Range("BU8").Value = 0.3682476943346509
$ws.Range("H9").Value = 3
$ws.Range("I9").Value = 8
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 8
$ws.Range("R9").Value = 7
$ws.Range("U9").Value = 8
$ws.Range("X9").Value = 8
$ws.Range("Y9").Value = 10
$ws.Range("AA9").Value = 6
$ws.Range("AD9").Value = 7
$ws.Range("AG9").Value = 6
$ws.Range("AJ9").Value = 6
$ws.Range("AM9").Value = 6
$ws.Range("AO9").Value = 26
$ws.Range("AP9").Value = 174
$ws.Range("AQ9").Value = 247
$ws.Range("AR9").Value = 18
$ws.Range("AS9").Value = 145
$ws.Range("AT9").Value = 246
$ws.Range("AU9").Value = 18
$ws.Range("AV9").Value = 130
$ws.Range("AW9").Value = 250
$ws.Range("AX9").Value = 12
$ws.Range("AY9").Value = 67
$ws.Range("AZ9").Value = 254
$ws.Range("BA9").Value = 13
$ws.Range("BB9").Value = 88
$ws.Range("BC9").Value = 255
$ws.Range("BD9").Value = 10
$ws.Range("BE9").Value = 101
$ws.Range("BF9").Value = 250
$ws.Range("BG9").Value = 0.2767811898246681
$ws.Range("BH9").Value = 0.4856298763940123
$ws.Range("BI9").Value = 0.7780041535708658
$ws.Range("BJ9").Value = 0.2414814814814815
$ws.Range("BK9").Value = 0.4562394127611517
$ws.Range("BL9").Value = 0.8539388543338676
$ws.Range("BM9").Value = 0.00206756368767997
$ws.Range("BN9").Value = 0.182918138489462
$ws.Range("BO9").Value = 0.35608683236975
$ws.Range("BP9").Value = 0.005068865814957325
$ws.Range("BQ9").Value = 0.1529829903340158
$ws.Range("BR9").Value = 0.3487574031177041
$ws.Range("BS9").Value = 0.000139366347672582
$ws.Range("BT9").Value = 0.04840579710144927
$ws.Range("BU9").Value = 0.4774866569626395
